$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (last-changed) date column C for all existing data rows (2-53)
$ws.Range("C2:C53").Value = 45206

# 2. Mark row 53 with an explicit (custom) row height, matching the rest of the sheet
$ws.Rows.Item(53).RowHeight = 15

# 3. Append a new data row (54) describing the new case
$ws.Range("A54").Value = "A 47874-2023"
$ws.Range("B54").Value = 45204
$ws.Range("C54").Value = 45206
$ws.Range("D54").Value = "SKÅNE LÄN"
$ws.Range("E54").Value = "SIMRISHAMN"
$ws.Range("F54").Value = "Övriga Aktiebolag"
$ws.Range("G54").Value = 1.4
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 0
$ws.Range("N54").Value = 0
$ws.Range("O54").Value = 0
$ws.Range("P54").Value = 0
$ws.Range("Q54").Value = 0

# Match the date formatting used by the rest of columns B/C
$ws.Range("B54").NumberFormat = $ws.Range("B53").NumberFormat
$ws.Range("C54").NumberFormat = $ws.Range("C53").NumberFormat

# R column keeps the wrap-text style used throughout the sheet (left empty, no species listed yet)
$ws.Range("R54").Style = $ws.Range("R53").Style
$ws.Range("R54").WrapText = $true
